$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 36) with the latest bitcoin buy entry, matching the
# plain-text date format used by the other recently-added rows.
$row = 36

# Leading apostrophe forces text entry (so the MM/DD/YYYY date string isn't
# auto-converted into a date serial number), then reset the style back to
# the workbook default so no number-format style gets attached to the cell
# (matching the plain, unstyled text cells used for the other recent rows).
$ws.Cells.Item($row, 1).Value = "'07/30/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.000422829999999999
$ws.Cells.Item($row, 3).Value = 118250.8336683776
$ws.Cells.Item($row, 4).Value = 50
